$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: roll the report forward one week.
#   "Volume 33   Number  4"                                  -> "...  5"
#   "Report Covering the Week  1/19/2026  Through  1/25/2026"
#     -> "Report Covering the Week  1/26/2026  Through  2/1/2026"
# Replace only the characters that actually changed (via Characters()) so
# the rest of each rich-text cell is left alone.
# ---------------------------------------------------------------------------
$cellA8 = $ws.Range("A8")
$a8Text = $cellA8.Text
$cellA8.Characters($a8Text.Length, 1).Text = "5"

$cellC9 = $ws.Range("C9")
$c9Text = $cellC9.Text
$weekStart = $c9Text.IndexOf("1/19/2026") + 1
$cellC9.Characters($weekStart, 9).Text = "1/26/2026"
$weekEnd = $c9Text.IndexOf("1/25/2026") + 1
$cellC9.Characters($weekEnd, 9).Text = "2/1/2026"

# ---------------------------------------------------------------------------
# A handful of cells flip between a numeric count and the "no data"
# placeholder text ("0" / "***.*"). Writing a different *type* into a cell
# through COM also resets its number format, so each conversion below:
#   1. Puts the cell into the right content type and value, then
#   2. Re-applies the original look by copying formats from a cell that
#      already has the target look (count cell vs. placeholder-text cell).
# ---------------------------------------------------------------------------

# C17: numeric 1 -> placeholder text "0"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "0"
$ws.Range("C22").Copy()
$ws.Range("C17").PasteSpecial(-4122)

# D27: numeric 1 -> placeholder text "0"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C22").Copy()
$ws.Range("D27").PasteSpecial(-4122)

# E27: numeric -100 -> placeholder text "***.*"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("E22").Copy()
$ws.Range("E27").PasteSpecial(-4122)

# C25: placeholder text "0" -> numeric 1
$ws.Range("C18").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C25").Value = 1

# I25: placeholder text "0" -> numeric 1
$ws.Range("C18").Copy()
$ws.Range("I25").PasteSpecial(-4122)
$ws.Range("I25").Value = 1

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Remaining cells: same content type, new value (freshly collected counts /
# recomputed percentages for the week).
# ---------------------------------------------------------------------------

# Row 16 - Robbery
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 3
$ws.Range("K16").Value = 66.666666666666
$ws.Range("L16").Value = -54.545454545454
$ws.Range("N16").Value = -80

# Row 17 - Fel. Assault
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 3
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -76.923076923076
$ws.Range("J17").Value = 14
$ws.Range("K17").Value = -71.428571428571
$ws.Range("L17").Value = -55.555555555555
$ws.Range("M17").Value = -33.333333333333
$ws.Range("N17").Value = -60

# Row 18 - Burglary
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 12
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = -40.740740740740
$ws.Range("I18").Value = 18
$ws.Range("J18").Value = 29
$ws.Range("K18").Value = -37.931034482758
$ws.Range("L18").Value = 5.882352941176
$ws.Range("M18").Value = -28
$ws.Range("N18").Value = -81.052631578947

# Row 19 - Gr. Larceny
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 0
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 47.058823529411
$ws.Range("I19").Value = 53
$ws.Range("J19").Value = 38
$ws.Range("K19").Value = 39.473684210526
$ws.Range("L19").Value = 55.882352941176
$ws.Range("M19").Value = 39.473684210526
$ws.Range("N19").Value = 20.454545454545

# Row 20 - G.L.A.
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 41.666666666666
$ws.Range("I20").Value = 17
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = -32
$ws.Range("M20").Value = 70
$ws.Range("N20").Value = -94.389438943894

# Row 21 - TOTAL
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -43.333333333333
$ws.Range("F21").Value = 90
$ws.Range("G21").Value = 88
$ws.Range("H21").Value = 2.272727272727
$ws.Range("I21").Value = 97
$ws.Range("J21").Value = 101
$ws.Range("K21").Value = -3.960396039603
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 8.988764044943
$ws.Range("N21").Value = -79.707112970711

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = -53.846153846153
$ws.Range("F24").Value = 37
$ws.Range("G24").Value = 39
$ws.Range("H24").Value = -5.128205128205
$ws.Range("I24").Value = 46
$ws.Range("J24").Value = 42
$ws.Range("K24").Value = 9.523809523809
$ws.Range("L24").Value = -19.298245614035
$ws.Range("M24").Value = -9.803921568627

# Row 25 - Retail Theft
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = -80
$ws.Range("J25").Value = 5
$ws.Range("K25").Value = -80
$ws.Range("L25").Value = -87.5

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 2
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 10
$ws.Range("G26").Value = 9
$ws.Range("H26").Value = 11.111111111111
$ws.Range("I26").Value = 11
$ws.Range("J26").Value = 10
$ws.Range("K26").Value = 10
$ws.Range("L26").Value = -15.384615384615
$ws.Range("M26").Value = -26.666666666666
